# Auto-generated Excel COM-interop script
# Applies numeric updates to the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets
# as described by the source diff (currentAveragePrice / profit columns).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 226.77777
$ws.Range("I53").Value = 214.42857
$ws.Range("K53").Value = 214.42857
$ws.Range("M53").Value = 422.57143
$ws.Range("H62").Value = 3507.3333
$ws.Range("I62").Value = 3209.3
$ws.Range("J62").Value = 4997.5
$ws.Range("K62").Value = 3209.3
$ws.Range("L62").Value = 4997.5
$ws.Range("M62").Value = -2585.3
$ws.Range("N62").Value = -6245.5
$ws.Range("H65").Value = 3507.3333
$ws.Range("I65").Value = 3209.3
$ws.Range("J65").Value = 4997.5
$ws.Range("K65").Value = 16046.5
$ws.Range("L65").Value = 24987.5
$ws.Range("M65").Value = -12926.5
$ws.Range("N65").Value = -31227.5
$ws.Range("H69").Value = 5000
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 5000
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 15000
$ws.Range("M69").ClearContents()
$ws.Range("N69").Value = -16748
$ws.Range("H72").Value = 5000
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 5000
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 45000
$ws.Range("M72").ClearContents()
$ws.Range("N72").Value = -53736
$ws.Range("H97").Value = 3875
$ws.Range("J97").Value = 3875
$ws.Range("L97").Value = 11625
$ws.Range("N97").Value = -12617
$ws.Range("H132").Value = 7904.2856
$ws.Range("I132").Value = 7904.2856
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 23712.8568
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -21182.8568
$ws.Range("N132").ClearContents()
$ws.Range("H138").Value = 3908.0613
$ws.Range("I138").Value = 3689
$ws.Range("J138").Value = 4004.7058
$ws.Range("K138").Value = 11067
$ws.Range("L138").Value = 12014.1174
$ws.Range("M138").Value = -5927
$ws.Range("N138").Value = -22294.1174
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1210836.6
$ws.Range("I32").Value = 1355972.9
$ws.Range("K32").Value = 1355972.9
$ws.Range("M32").Value = -1355685.9
$ws.Range("H61").Value = 2943237.5
$ws.Range("I61").Value = 1769.4667
$ws.Range("K61").Value = 1769.4667
$ws.Range("M61").Value = -1557.4667
$ws.Range("H74").Value = 875865.0600000001
$ws.Range("I74").Value = 1015473.44
$ws.Range("K74").Value = 1015473.44
$ws.Range("M74").Value = -1014599.44
$ws.Range("H77").Value = 875865.0600000001
$ws.Range("I77").Value = 1015473.44
$ws.Range("K77").Value = 5077367.199999999
$ws.Range("M77").Value = -5072999.199999999
$ws.Range("H97").Value = 691.44446
$ws.Range("I97").Value = 691.44446
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 691.44446
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -195.44446
$ws.Range("N97").ClearContents()
$ws.Range("H129").Value = 90000
$ws.Range("J129").Value = 90000
$ws.Range("L129").Value = 90000
$ws.Range("N129").Value = -100000
$ws.Range("H131").Value = 50000
$ws.Range("J131").Value = 50000
$ws.Range("L131").Value = 50000
$ws.Range("N131").Value = -60080
$ws.Range("H136").Value = 2943237.5
$ws.Range("I136").Value = 1769.4667
$ws.Range("K136").Value = 5308.4001
$ws.Range("M136").Value = -2758.4001
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 53731.188
$ws.Range("I20").Value = 61180.48
$ws.Range("K20").Value = 61180.48
$ws.Range("M20").Value = -60933.48
$ws.Range("H94").Value = 5358.3335
$ws.Range("J94").Value = 1010
$ws.Range("L94").Value = 1010
$ws.Range("N94").Value = -1912
$ws.Range("H107").Value = 1389.9459
$ws.Range("I107").Value = 1064.7273
$ws.Range("J107").Value = 1866.9333
$ws.Range("K107").Value = 1064.7273
$ws.Range("L107").Value = 1866.9333
$ws.Range("M107").Value = 855.2727
$ws.Range("N107").Value = -5706.9333
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8236882
$ws.Range("I31").Value = 2858499.8
$ws.Range("K31").Value = 2858499.8
$ws.Range("M31").Value = -2858204.8
$ws.Range("H34").Value = 8236882
$ws.Range("I34").Value = 2858499.8
$ws.Range("K34").Value = 2858499.8
$ws.Range("M34").Value = -2858297.8
$ws.Range("H107").Value = 618.1818
$ws.Range("I107").Value = 406.35715
$ws.Range("J107").Value = 988.875
$ws.Range("K107").Value = 406.35715
$ws.Range("L107").Value = 988.875
$ws.Range("M107").Value = 1513.64285
$ws.Range("N107").Value = -4828.875
$ws.Range("H134").Value = 2530.4062
$ws.Range("I134").Value = 2599.1724
$ws.Range("K134").Value = 7797.5172
$ws.Range("M134").Value = -5262.5172
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 9811.546
$ws.Range("I3").Value = 3785.4
$ws.Range("K3").Value = 11356.2
$ws.Range("M3").Value = -11244.2
$ws.Range("H37").Value = 46315.79
$ws.Range("J37").Value = 46315.79
$ws.Range("L37").Value = 138947.37
$ws.Range("N37").Value = -139171.37
$ws.Range("H101").Value = 6679150.5
$ws.Range("J101").Value = 6679150.5
$ws.Range("L101").Value = 20037451.5
$ws.Range("N101").Value = -20042319.5
$ws.Range("H131").Value = 2437.244
$ws.Range("I131").Value = 576.6667
$ws.Range("J131").Value = 3510.6538
$ws.Range("K131").Value = 1730.0001
$ws.Range("L131").Value = 10531.9614
$ws.Range("M131").Value = 3309.9999
$ws.Range("N131").Value = -20611.9614
$ws.Range("H133").Value = 9872.786
$ws.Range("I133").Value = 5304.8335
$ws.Range("J133").Value = 13298.75
$ws.Range("K133").Value = 15914.5005
$ws.Range("L133").Value = 39896.25
$ws.Range("M133").Value = -10854.5005
$ws.Range("N133").Value = -50016.25
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 18318.096
$ws.Range("I70").Value = 35781.168
$ws.Range("K70").Value = 35781.168
$ws.Range("M70").Value = -35511.168
$ws.Range("H73").Value = 18318.096
$ws.Range("I73").Value = 35781.168
$ws.Range("K73").Value = 35781.168
$ws.Range("M73").Value = -34845.168
$ws.Range("H80").Value = 1832.6666
$ws.Range("I80").Value = 1499
$ws.Range("K80").Value = 1499
$ws.Range("M80").Value = -501
$ws.Range("H83").Value = 1832.6666
$ws.Range("I83").Value = 1499
$ws.Range("K83").Value = 7495
$ws.Range("M83").Value = -2503
$ws.Range("H97").Value = 1261.5385
$ws.Range("I97").Value = 1056.7368
$ws.Range("K97").Value = 1056.7368
$ws.Range("M97").Value = -560.7367999999999
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H88").Value = 39996.668
$ws.Range("J88").Value = 39996.668
$ws.Range("L88").Value = 39996.668
$ws.Range("N88").Value = -40852.668
$ws.Range("H91").Value = 39996.668
$ws.Range("J91").Value = 39996.668
$ws.Range("L91").Value = 39996.668
$ws.Range("N91").Value = -42960.668
$ws.Range("H122").Value = 4274.3335
$ws.Range("J122").Value = 4511.4287
$ws.Range("L122").Value = 13534.2861
$ws.Range("N122").Value = -18434.2861
$ws.Range("H130").Value = 85000
$ws.Range("J130").Value = 85000
$ws.Range("L130").Value = 85000
$ws.Range("N130").Value = -95040
$ws.Range("H132").Value = 3208852
$ws.Range("I132").Value = 5956336.5
$ws.Range("J132").Value = 3453.4167
$ws.Range("K132").Value = 17869009.5
$ws.Range("L132").Value = 10360.2501
$ws.Range("M132").Value = -17866479.5
$ws.Range("N132").Value = -15420.2501
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3833.3333
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()
$ws.Range("H84").Value = 3833.3333
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()
$ws.Range("H122").Value = 44608.652
$ws.Range("I122").Value = 1261.2609
$ws.Range("J122").Value = 376938.66
$ws.Range("K122").Value = 3783.7827
$ws.Range("L122").Value = 1130815.98
$ws.Range("N122").Value = -1133715.98
